$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''64.445.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +0.26%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.141.44'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -0.30%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.04%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''608.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +0.30%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''144.02'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -2.25%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  -0.09%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''3.140.24'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -0.21%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '''  +1.26%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '''  -0.36%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  -2.26%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '''  -0.93%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  +1.75%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''35.47'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -0.52%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''3.658.26'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -0.38%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = '''  +2.47%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''64.510.79'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +0.30%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''3.141.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -0.51%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''6.87'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -0.86%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''477.57'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -0.73%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = '''  +0.58%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.719'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +1.11%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''7.80'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +0.16%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''85.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +2.35%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''13.51'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -1.28%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''  +0.08%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''2.78'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -3.48%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''8.50'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +0.46%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''7.33'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +7.40%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''0.116'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +1.78%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = '''  -5.90%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  -0.07%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''26.75'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +1.99%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''2.65'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -3.26%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  +0.77%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''5.98'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +0.05%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''52.71'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -3.08%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.0₃0743'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +1.58%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''3.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +2.20%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''450.02'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -0.80%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.0395'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -0.27%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''  -0.35%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  -1.31%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''2.890.29'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +1.43%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.263'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -1.13%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''2.24'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -1.05%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''2.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +5.03%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''26.42'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -0.70%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E50").Value = '''  -0.14%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''120.83'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +0.89%  '
$ws.Range("E51").Style = "Normal"
